# Atualizado por script em 21-12-2023 02:45
#
# 1) Rows 16 and 17 (match data in columns F:V) were swapped - the
#    "Maccabi Tel Aviv v Maccabi Bnei Raina" match moved from row 16 to
#    row 17, and "Hapoel Haifa v SC Ashdod" moved from row 17 to row 16.
#    Column A (the running index) stays tied to the row number.
# 2) Four new match rows (67-70) were appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: swap the contents of row 16 and row 17 (columns F through V)
# ---------------------------------------------------------------------
$cols = 6..22   # F=6 ... V=22

foreach ($c in $cols) {
    $v16 = $ws.Cells.Item(16, $c).Value2
    $v17 = $ws.Cells.Item(17, $c).Value2
    $ws.Cells.Item(16, $c).Value2 = $v17
    $ws.Cells.Item(17, $c).Value2 = $v16
}

# ---------------------------------------------------------------------
# Step 2: append four new rows (67-70) with the same look & feel
# (number formats / fonts / borders) as the last existing row (66).
# ---------------------------------------------------------------------
$ws.Range("A66:V66").Copy()
$ws.Range("A67:V70").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{ Row=67; A=66; B="israel"; C="ligat-ha-al"; D="2023-2024"; E=45280.75;              F="Sakhnin";         G=1; H="Hapoel Haifa";     I=1; J=3.06; K="17/12/2023 19:43"; L=3.28; M="20/12/2023 17:56"; N=3.17; O="17/12/2023 19:43"; P=3.3;  Q="20/12/2023 17:56"; R=2.31; S="17/12/2023 19:43"; T=2.29; U="20/12/2023 17:56"; V="https://www.betexplorer.com/football/israel/ligat-ha-al/sakhnin-hapoel-haifa/4vbajJ8g/" }
    @{ Row=68; A=67; B="israel"; C="ligat-ha-al"; D="2023-2024"; E=45280.78125;           F="Netanya";         G=0; H="Beitar Jerusalem"; I=3; J=2.18; K="17/12/2023 19:12"; L=1.99; M="20/12/2023 18:44"; N=3.34; O="17/12/2023 19:12"; P=3.58; Q="20/12/2023 18:42"; R=3.16; S="17/12/2023 19:12"; T=3.8;  U="20/12/2023 18:42"; V="https://www.betexplorer.com/football/israel/ligat-ha-al/netanya-beitar-jerusalem/Is23kwOa/" }
    @{ Row=69; A=68; B="israel"; C="ligat-ha-al"; D="2023-2024"; E=45280.79166666666;     F="Maccabi Haifa";   G=4; H="SC Ashdod";         I=0; J=1.32; K="17/12/2023 19:43"; L=1.31; M="20/12/2023 18:55"; N=4.98; O="17/12/2023 19:43"; P=5.66; Q="20/12/2023 18:55"; R=7.68; S="17/12/2023 19:43"; T=8.84; U="20/12/2023 18:55"; V="https://www.betexplorer.com/football/israel/ligat-ha-al/maccabi-haifa-sc-ashdod/MF8fiagm/" }
    @{ Row=70; A=69; B="israel"; C="ligat-ha-al"; D="2023-2024"; E=45280.8125;            F="H. Beer Sheva";   G=0; H="Maccabi Tel Aviv";  I=1; J=4.71; K="17/12/2023 19:12"; L=4.19; M="20/12/2023 19:28"; N=3.74; O="17/12/2023 19:12"; P=3.68; Q="20/12/2023 19:28"; R=1.67; S="17/12/2023 19:12"; T=1.86; U="20/12/2023 19:28"; V="https://www.betexplorer.com/football/israel/ligat-ha-al/h-beer-sheva-maccabi-tel-aviv/C417lcw6/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = $r.B
    $ws.Cells.Item($row, 3).Value2  = $r.C
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
    $ws.Cells.Item($row, 21).Value2 = $r.U
    $ws.Cells.Item($row, 22).Value2 = $r.V
}

Write-Output "done"
